$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A30").Value = "24/10/2025"
$ws.Range("B30").Value = "Sarmiento Junin"
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = "Rosario"
$ws.Range("F30").Value = "D"
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 1
$ws.Range("K30").Value = 0.09
$ws.Range("L30").Value = 1.28
$ws.Range("M30").Value = 3
$ws.Range("N30").Value = 7
$ws.Range("O30").Value = 2
$ws.Range("P30").Value = 5
